$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dutch")

# Insert a new row before row 87 (pushes existing rows 87-102 down to 88-103),
# inheriting formatting from the row below (same as Excel's default insert
# behaviour / "Format Same As Below" here since the block below shares style).
$ws.Rows.Item(87).Insert()

# New "alumni" / podcast news-item entry: misdadiger / misdadiger (criminal, more felonious)
$ws.Range("A87").Value = 4
$ws.Range("B87").Value = "misdadiger"
$ws.Range("C87").Value = 1
$ws.Range("D87").Value = "criminal (noun)"
$ws.Range("E87").Value = "."
$ws.Range("F87").Value = "misdadiger"
$ws.Range("G87").Value = 2
$ws.Range("H87").Value = "more felonious"

# Tidy up the sheet view: clear the frozen/scrolled top-left cell and the
# lingering selection on H87 left over from editing.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("A1").Select()
